# Fruta / hortaliza, semanal
# Insert a new daily price record as row 85 (Palta, Hass, Primera) for the
# "Agrícola del Norte S.A. de Arica" market; every subsequent row shifts
# down by one (85->86, 86->87, ..., 141->142).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row above the current row 85; Excel shifts rows 85..141
# down to 86..142, carrying their content and formatting along.
$ws.Rows.Item(85).Insert()

# Populate the newly inserted row 85 with the new record's values.
$ws.Cells.Item(85, 1).Value = 1
$ws.Cells.Item(85, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(85, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(85, 4).Value = 44879
$ws.Cells.Item(85, 5).Value = 15
$ws.Cells.Item(85, 6).Value = "Fruta"
$ws.Cells.Item(85, 7).Value = 100106
$ws.Cells.Item(85, 8).Value = "Oleaginosos"
$ws.Cells.Item(85, 9).Value = 100106002
$ws.Cells.Item(85, 10).Value = "Palta"
$ws.Cells.Item(85, 11).Value = "Hass"
$ws.Cells.Item(85, 12).Value = "Primera"
$ws.Cells.Item(85, 13).Value = 416
$ws.Cells.Item(85, 14).Value = 27000
$ws.Cells.Item(85, 15).Value = 28000
$ws.Cells.Item(85, 16).Value = 27500
$ws.Cells.Item(85, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(85, 18).Value = "Perú"
$ws.Cells.Item(85, 19).Value = 2750
$ws.Cells.Item(85, 20).Value = 10
